$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.028784482558156
$ws.Range("D2").Value = 1.031106560225712
$ws.Range("E2").Value = 1.037373818706157
$ws.Range("F2").Value = 1.045197633122812
$ws.Range("I2").Value = 1.030645842577517
$ws.Range("J2").Value = 1.033934879726482
$ws.Range("K2").Value = 1.03391563697029
$ws.Range("L2").Value = 1.040164886640484
$ws.Range("M2").Value = 1.047966548433424
$ws.Range("N2").Value = 1.005712725503983
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.029788170279179
$ws.Range("D3").Value = 1.031822661745746
$ws.Range("E3").Value = 1.038300522135754
$ws.Range("F3").Value = 1.046291838035137
$ws.Range("I3").Value = 1.030808273457226
$ws.Range("J3").Value = 1.034579045482416
$ws.Range("K3").Value = 1.034440480933844
$ws.Range("L3").Value = 1.040901074917712
$ws.Range("M3").Value = 1.048871408493035
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.030437600828379
$ws.Range("D4").Value = 1.032285477444552
$ws.Range("E4").Value = 1.038900590218199
$ws.Range("F4").Value = 1.047000598492137
$ws.Range("I4").Value = 1.030911346223459
$ws.Range("J4").Value = 1.034995266394261
$ws.Range("K4").Value = 1.034778872226354
$ws.Range("L4").Value = 1.041377241684063
$ws.Range("M4").Value = 1.049457058461462
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.030710615713436
$ws.Range("D5").Value = 1.032479912409638
$ws.Range("E5").Value = 1.039152960796598
$ws.Range("F5").Value = 1.047298737200053
$ws.Range("I5").Value = 1.030954191198927
$ws.Range("J5").Value = 1.035170102035964
$ws.Range("K5").Value = 1.034920839546412
$ws.Range("L5").Value = 1.041577374598863
$ws.Range("M5").Value = 1.049703299845838
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.030756455769145
$ws.Range("D6").Value = 1.032512551078431
$ws.Range("E6").Value = 1.039195340884218
$ws.Range("F6").Value = 1.047348806294415
$ws.Range("I6").Value = 1.030961356493581
$ws.Range("J6").Value = 1.0351994492692
$ws.Range("K6").Value = 1.034944659314749
$ws.Range("L6").Value = 1.041610974987312
$ws.Range("M6").Value = 1.049744646842831
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.030441248889518
$ws.Range("D7").Value = 1.032288076017132
$ws.Range("E7").Value = 1.038903962006546
$ws.Range("F7").Value = 1.047004581546031
$ws.Range("I7").Value = 1.030911920633764
$ws.Range("J7").Value = 1.034997603121808
$ws.Range("K7").Value = 1.034780770350541
$ws.Range("L7").Value = 1.041379916057648
$ws.Range("M7").Value = 1.049460348617487
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.029123687990597
$ws.Range("D8").Value = 1.031348683273746
$ws.Range("E8").Value = 1.037686913194885
$ws.Range("F8").Value = 1.045567271919691
$ws.Range("I8").Value = 1.030701156989109
$ws.Range("J8").Value = 1.034152701891751
$ws.Range("K8").Value = 1.03409326211532
$ws.Range("L8").Value = 1.04041372519336
$ws.Range("M8").Value = 1.048272319757126
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.026801818607665
$ws.Range("D9").Value = 1.029689190448959
$ws.Range("E9").Value = 1.035545635706724
$ws.Range("F9").Value = 1.043040215130841
$ws.Range("I9").Value = 1.030314237353928
$ws.Range("J9").Value = 1.032659328914761
$ws.Range("K9").Value = 1.032872493271132
$ws.Range("L9").Value = 1.038709696391907
$ws.Range("M9").Value = 1.046179996570818
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.025253814436939
$ws.Range("D10").Value = 1.028580126586299
$ws.Range("E10").Value = 1.034120389714618
$ws.Range("F10").Value = 1.041359350138021
$ws.Range("I10").Value = 1.030045886450925
$ws.Range("J10").Value = 1.031660721950384
$ws.Range("K10").Value = 1.032052446663378
$ws.Range("L10").Value = 1.037572717384435
$ws.Range("M10").Value = 1.044785909037276
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.02458349217619
$ws.Range("D11").Value = 1.028099253690401
$ws.Range("E11").Value = 1.033503789475113
$ws.Range("F11").Value = 1.040632433129845
$ws.Range("I11").Value = 1.0299272267649
$ws.Range("J11").Value = 1.031227603061432
$ws.Range("K11").Value = 1.031695895887648
$ws.Range("L11").Value = 1.037080171633253
$ws.Range("M11").Value = 1.044182448012637
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.024334500767585
$ws.Range("D12").Value = 1.027920540753367
$ws.Range("E12").Value = 1.033274838493953
$ws.Range("F12").Value = 1.040362560637946
$ws.Range("I12").Value = 1.029882781895308
$ws.Range("J12").Value = 1.031066616368232
$ws.Range("K12").Value = 1.031563237590666
$ws.Range("L12").Value = 1.036897184641552
$ws.Range("M12").Value = 1.04395832435416
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.024387910444245
$ws.Range("D13").Value = 1.027958879583721
$ws.Range("E13").Value = 1.033323945545775
$ws.Range("F13").Value = 1.040420443005336
$ws.Range("I13").Value = 1.029892332177841
$ws.Range("J13").Value = 1.031101153397776
$ws.Range("K13").Value = 1.031591703166298
$ws.Range("L13").Value = 1.036936437484557
$ws.Range("M13").Value = 1.044006398337861
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.024562910547347
$ws.Range("D14").Value = 1.02808448316252
$ws.Range("E14").Value = 1.033484862643568
$ws.Range("F14").Value = 1.040610122591802
$ws.Range("I14").Value = 1.029923560470435
$ws.Range("J14").Value = 1.031214298032477
$ws.Range("K14").Value = 1.031684934784614
$ws.Range("L14").Value = 1.037065046553731
$ws.Range("M14").Value = 1.044163921294784
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.024670733375118
$ws.Range("D15").Value = 1.028161859068883
$ws.Range("E15").Value = 1.033584019835006
$ws.Range("F15").Value = 1.04072700858945
$ws.Range("I15").Value = 1.029942752320873
$ws.Range("J15").Value = 1.031283996003651
$ws.Range("K15").Value = 1.031742348800371
$ws.Range("L15").Value = 1.037144282410615
$ws.Range("M15").Value = 1.044260980193173
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.025298300970699
$ws.Range("D16").Value = 1.028612027115222
$ws.Range("E16").Value = 1.034161322882122
$ws.Range("F16").Value = 1.041407612398604
$ws.Range("I16").Value = 1.030053709675769
$ws.Range("J16").Value = 1.031689451587864
$ws.Range("K16").Value = 1.032076078951206
$ws.Range("L16").Value = 1.037605401292372
$ws.Range("M16").Value = 1.044825962752997
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.025691950566819
$ws.Range("D17").Value = 1.028894234723932
$ws.Range("E17").Value = 1.03452359499138
$ws.Range("F17").Value = 1.041834780692859
$ws.Range("I17").Value = 1.030122651511215
$ws.Range("J17").Value = 1.031943591920998
$ws.Range("K17").Value = 1.03228502725792
$ws.Range("L17").Value = 1.037894588696134
$ws.Range("M17").Value = 1.045180412008273
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.025921557180136
$ws.Range("D18").Value = 1.02905877974382
$ws.Range("E18").Value = 1.034734954313334
$ws.Range("F18").Value = 1.042084028378346
$ws.Range("I18").Value = 1.030162626539124
$ws.Range("J18").Value = 1.032091758614587
$ws.Range("K18").Value = 1.03240676185712
$ws.Range("L18").Value = 1.038063244948772
$ws.Range("M18").Value = 1.045387174383569
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.025999846669454
$ws.Range("D19").Value = 1.029114874814111
$ws.Range("E19").Value = 1.034807031204787
$ws.Range("F19").Value = 1.042169030253943
$ws.Range("I19").Value = 1.030176216677734
$ws.Range("J19").Value = 1.032142267900379
$ws.Range("K19").Value = 1.032448246203135
$ws.Range("L19").Value = 1.038120748650061
$ws.Range("M19").Value = 1.045457678072471
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.025649715952317
$ws.Range("D20").Value = 1.028863962914275
$ws.Range("E20").Value = 1.034484721226317
$ws.Range("F20").Value = 1.041788940509296
$ws.Range("I20").Value = 1.030115279274865
$ws.Range("J20").Value = 1.031916332208654
$ws.Range("K20").Value = 1.032262623704202
$ws.Range("L20").Value = 1.037863563883154
$ws.Range("M20").Value = 1.045142381080423
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.024511377485183
$ws.Range("D21").Value = 1.028047498663144
$ws.Range("E21").Value = 1.033437474306325
$ws.Range("F21").Value = 1.040554262900752
$ws.Range("I21").Value = 1.029914374709253
$ws.Range("J21").Value = 1.031180982709297
$ws.Range("K21").Value = 1.031657486449732
$ws.Range("L21").Value = 1.037027175311549
$ws.Range("M21").Value = 1.04411753391663
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.023795636523233
$ws.Range("D22").Value = 1.027533604085711
$ws.Range("E22").Value = 1.032779502162881
$ws.Range("F22").Value = 1.039778764027915
$ws.Range("I22").Value = 1.029785921168039
$ws.Range("J22").Value = 1.030718020130533
$ws.Range("K22").Value = 1.031275743215042
$ws.Range("L22").Value = 1.036501110612731
$ws.Range("M22").Value = 1.043473337688232
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.024175066449704
$ws.Range("D23").Value = 1.027806081244408
$ws.Range("E23").Value = 1.033128260567862
$ws.Range("F23").Value = 1.040189795473325
$ws.Range("I23").Value = 1.029854219204063
$ws.Range("J23").Value = 1.030963503896917
$ws.Range("K23").Value = 1.031478232633884
$ws.Range("L23").Value = 1.036780005708352
$ws.Range("M23").Value = 1.043814822450355
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.025668799964337
$ws.Range("D24").Value = 1.028877641633324
$ws.Range("E24").Value = 1.03450228644634
$ws.Range("F24").Value = 1.041809653444177
$ws.Range("I24").Value = 1.030118611205218
$ws.Range("J24").Value = 1.03192864991263
$ws.Range("K24").Value = 1.03227274734271
$ws.Range("L24").Value = 1.037877582728211
$ws.Range("M24").Value = 1.04515956556363
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.027402094700151
$ws.Range("D25").Value = 1.030118695850171
$ws.Range("E25").Value = 1.036098809767904
$ws.Range("F25").Value = 1.043692845305579
$ws.Range("I25").Value = 1.030416101758365
$ws.Range("J25").Value = 1.033045937233087
$ws.Range("K25").Value = 1.03318918755799
$ws.Range("L25").Value = 1.039150400273363
$ws.Range("M25").Value = 1.046720774481905
